# Weekly fruit/vegetable price update: insert a new data row at row 369
# (most-recent-week-first ordering), pushing the existing rows 369-379
# down to 370-380. The new row carries a fresh observation; all other
# rows keep their original values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 369, shifting rows 369:379 -> 370:380.
$ws.Rows(369).Insert()

# Populate the newly inserted row 369 with the new weekly observation.
$ws.Range("A369").Value = 6
$ws.Range("B369").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C369").Value = "Metropolitana"
$ws.Range("D369").Value = 45239
$ws.Range("E369").Value = 13
$ws.Range("F369").Value = 100112001
$ws.Range("G369").Value = "Berenjena"
$ws.Range("H369").Value = "Sin especificar"
$ws.Range("I369").Value = "Primera"
$ws.Range("J369").Value = 400
$ws.Range("K369").Value = 9000
$ws.Range("L369").Value = 10000
$ws.Range("M369").Value = 9425
$ws.Range("N369").Value = "`$/caja 50 unidades"
$ws.Range("O369").Value = "Región de Arica y Parinacota"
$ws.Range("P369").Value = 188
$ws.Range("Q369").Value = 50
$ws.Range("R369").Value = "Hortaliza"
